# Updates the "Year Figures" table for 2024 (rows 191-200 on Sheet1) with
# the latest December snapshot of points/chips/winnings/takehome/pers_personid,
# including a couple of swapped rankings between people.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: A=Yr  B=Person  C=SRank  D=Points  E=Bonus  F=PointsBonus
#          G=Chips  H=Winnings  I=Takehome  J=PersStatus  K=pers_personid

$rows = @(
    @{ Row = 191; Person = "Richard";  SRank = 1;  Points = 63; Bonus = 0; PointsBonus = 63; Chips = 206300; Winnings = 200; Takehome = 100;  PersId = 366 }
    @{ Row = 192; Person = "Mark";     SRank = 2;  Points = 47; Bonus = 0; PointsBonus = 47; Chips = 140750; Winnings = 180; Takehome = 70;   PersId = 361 }
    @{ Row = 193; Person = "Andy";     SRank = 3;  Points = 40; Bonus = 0; PointsBonus = 40; Chips = 138000; Winnings = 130; Takehome = 30;   PersId = 349 }
    @{ Row = 194; Person = "Anthony";  SRank = 4;  Points = 32; Bonus = 0; PointsBonus = 32; Chips = 108750; Winnings = 40;  Takehome = -40;  PersId = 350 }
    @{ Row = 195; Person = "Matt";     SRank = 5;  Points = 30; Bonus = 0; PointsBonus = 30; Chips = 116550; Winnings = 40;  Takehome = -70;  PersId = 362 }
    @{ Row = 196; Person = "Prashant"; SRank = 6;  Points = 28; Bonus = 0; PointsBonus = 28; Chips = 97650;  Winnings = 130; Takehome = 40;   PersId = 365 }
    @{ Row = 197; Person = "Pepe";     SRank = 7;  Points = 24; Bonus = 0; PointsBonus = 24; Chips = 77350;  Winnings = 40;  Takehome = -40;  PersId = 364 }
    @{ Row = 198; Person = "Jon";      SRank = 8;  Points = 23; Bonus = 0; PointsBonus = 23; Chips = 76650;  Winnings = 20;  Takehome = -90;  PersId = 357 }
    @{ Row = 199; Person = "Maisy";    SRank = 9;  Points = 20; Bonus = 0; PointsBonus = 20; Chips = 74950;  Winnings = 50;  Takehome = -40;  PersId = 360 }
    @{ Row = 200; Person = "Alex";     SRank = 10; Points = 19; Bonus = 1; PointsBonus = 20; Chips = 69950;  Winnings = 80;  Takehome = 40;   PersId = 348 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 2).Value = $r.Person
    $ws.Cells.Item($rowNum, 3).Value = $r.SRank
    $ws.Cells.Item($rowNum, 4).Value = $r.Points
    $ws.Cells.Item($rowNum, 5).Value = $r.Bonus
    $ws.Cells.Item($rowNum, 6).Value = $r.PointsBonus
    $ws.Cells.Item($rowNum, 7).Value = $r.Chips
    $ws.Cells.Item($rowNum, 8).Value = $r.Winnings
    $ws.Cells.Item($rowNum, 9).Value = $r.Takehome
    $ws.Cells.Item($rowNum, 11).Value = $r.PersId
}
